$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete "ECs" sending-cluster rows (old rows 8-10 stay removed;
# old rows 2-4 (ECs sender) are overwritten below with the new FAPs sender values,
# so only the trailing rows 8:10 need to be deleted).
$ws.Range("A8:T10").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.900731333333333
$ws.Range("H2").Value = 8.702194
$ws.Range("I2").Value = 0.8130494232775288
$ws.Range("J2").Value = 0.8130494232775289
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.030956000000001
$ws.Range("N2").Value = 18.092868
$ws.Range("O2").Value = 0.364814105361131
$ws.Range("P2").Value = 0.3648141053611309
$ws.Range("Q2").Value = 17.49418303915467
$ws.Range("R2").Value = 157.447647352392
$ws.Range("S2").Value = 0.2966118979673751
$ws.Range("T2").Value = 0.2966118979673751

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.900731333333333
$ws.Range("H3").Value = 8.702194
$ws.Range("I3").Value = 0.8130494232775288
$ws.Range("J3").Value = 0.8130494232775289
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.789877333333333
$ws.Range("N3").Value = 20.369632
$ws.Range("O3").Value = 0.4107214552505144
$ws.Range("P3").Value = 0.4107214552505143
$ws.Range("Q3").Value = 19.69560993028978
$ws.Range("R3").Value = 177.260489372608
$ws.Range("S3").Value = 0.3339368423191381
$ws.Range("T3").Value = 0.3339368423191381

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.900731333333333
$ws.Range("H4").Value = 8.702194
$ws.Range("I4").Value = 0.8130494232775288
$ws.Range("J4").Value = 0.8130494232775289
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.710753333333333
$ws.Range("N4").Value = 11.13226
$ws.Range("O4").Value = 0.2244644393883547
$ws.Range("P4").Value = 0.2244644393883547
$ws.Range("Q4").Value = 10.76389846427111
$ws.Range("R4").Value = 96.87508617844
$ws.Range("S4").Value = 0.1825006829910156
$ws.Range("T4").Value = 0.1825006829910156

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6669870000000001
$ws.Range("H5").Value = 2.000961
$ws.Range("I5").Value = 0.1869505767224711
$ws.Range("J5").Value = 0.1869505767224711
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.030956000000001
$ws.Range("N5").Value = 18.092868
$ws.Range("O5").Value = 0.364814105361131
$ws.Range("P5").Value = 0.3648141053611309
$ws.Range("Q5").Value = 4.022569249572001
$ws.Range("R5").Value = 36.20312324614801
$ws.Range("S5").Value = 0.06820220739375576
$ws.Range("T5").Value = 0.06820220739375575

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6669870000000001
$ws.Range("H6").Value = 2.000961
$ws.Range("I6").Value = 0.1869505767224711
$ws.Range("J6").Value = 0.1869505767224711
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.789877333333333
$ws.Range("N6").Value = 20.369632
$ws.Range("O6").Value = 0.4107214552505144
$ws.Range("P6").Value = 0.4107214552505143
$ws.Range("Q6").Value = 4.528759912928001
$ws.Range("R6").Value = 40.758839216352
$ws.Range("S6").Value = 0.07678461293137627
$ws.Range("T6").Value = 0.07678461293137624

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6669870000000001
$ws.Range("H7").Value = 2.000961
$ws.Range("I7").Value = 0.1869505767224711
$ws.Range("J7").Value = 0.1869505767224711
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.710753333333333
$ws.Range("N7").Value = 11.13226
$ws.Range("O7").Value = 0.2244644393883547
$ws.Range("P7").Value = 0.2244644393883547
$ws.Range("Q7").Value = 2.47502423354
$ws.Range("R7").Value = 22.27521810186
$ws.Range("S7").Value = 0.04196375639733907
$ws.Range("T7").Value = 0.04196375639733906

